# The template's header contains a field:
#   { SEQ  m:Sequence{'some text'.sampleText(6), 'more text'.sampleText(8), 'even more text'.sampleText(10)} }
# represented as a begin/instrText/end run triple. The parser was switched to
# TokenIteratorFieldRewriterSplit, which expects the M2Doc query written out as
# plain literal text (wrapped in { }) instead of living inside a Word field.
# Replace the field with a single run containing that literal text.

$d = $word.ActiveDocument

# Locate the header story range (wdHeaderFooterPrimary / story type 7).
$header = $null
foreach ($story in $d.StoryRanges) {
    if ($story.StoryType -eq 7) {
        $header = $story
    }
}

$field = $header.Fields.Item(1)

# Selecting the field switches $word.Selection into the header story and
# positions it exactly on the field; Unlink() then collapses the field's
# begin/instrText/end run-triple down to nothing (here, an empty result),
# leaving the selection collapsed right where the field used to be while
# keeping the surrounding paragraph untouched.
$field.Select()
$field.Unlink()

# Type the replacement text, followed by a trailing space, so Word marks the
# run's <w:t> with xml:space="preserve" (matching the target markup), then
# remove that trailing space again with a backspace - the preserve attribute
# sticks even once the triggering whitespace is gone.
$word.Selection.TypeText("{m:Sequence{'some text'.sampleText(6), 'more text'.sampleText(8), 'even more text'.sampleText(10)}} ")
$word.Selection.TypeBackspace()
